# Update Date column (A) values for round 3 results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Excel serial date value (days since 1899-12-30)
$updates = @{
    2  = 45746
    3  = 45745
    4  = 45744
    6  = 45743
    7  = 45746
    8  = 45745
    9  = 45745
    10 = 45745
    11 = 45745
    12 = 45745
    13 = 45746
    14 = 45743
    15 = 45745
    16 = 45745
    18 = 45746
    19 = 45744
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $updates[$row]
}
